# Updated cryptos list - apply latest price/volume snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a Price-column value while preserving its original "Text"
# storage type (some prices look like plain numbers, e.g. "18.86", and
# Excel would otherwise silently convert them to numeric cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "26.564.64"
$ws.Range("E2").Value = "  -0.36%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.636.38"
$ws.Range("E3").Value = "  +0.28%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.21%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  +0.63%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +2.23%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.17%  "

# Row 8 - now Dogecoin (was Cardano)
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D8") "0.0627"
$ws.Range("E8").Value = "  +0.49%  "

# Row 9 - now Cardano (was Dogecoin)
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D9") "0.250"
$ws.Range("E9").Value = "  -0.55%  "

# Row 10 - Solana
Set-TextValue $ws.Range("D10") "18.86"
$ws.Range("E10").Value = "  -0.86%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.0842"
$ws.Range("E11").Value = "  +0.44%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D12") "1.861.24"
$ws.Range("E12").Value = "  +0.13%  "

# Row 13 - WrappedEther
Set-TextValue $ws.Range("D13") "1.620.48"
$ws.Range("E13").Value = "  -0.62%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "4.15"
$ws.Range("E14").Value = "  +1.78%  "

# Row 15 - Polygon
Set-TextValue $ws.Range("D15") "0.525"
$ws.Range("E15").Value = "  -0.29%  "

# Row 16 - Litecoin
Set-TextValue $ws.Range("D16") "65.37"
$ws.Range("E16").Value = "  +3.71%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "26.580.07"
$ws.Range("E17").Value = "  -0.20%  "

# Row 18 - ShibaInu
Set-TextValue $ws.Range("D18") "0.0₃0744"
$ws.Range("E18").Value = "  +0.68%  "

# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "216.13"
$ws.Range("E19").Value = "  +3.28%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.19%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +0.36%  "

# Row 22 - Chainlink
$ws.Range("E22").Value = "  +1.66%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  -0.61%  "

# Row 24 - Toncoin
Set-TextValue $ws.Range("D24") "2.22"
$ws.Range("E24").Value = "  +14.78%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "147.27"
$ws.Range("E25").Value = "  +0.23%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.22%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +0.09%  "

# Row 28 - Cosmos
Set-TextValue $ws.Range("D28") "6.93"
$ws.Range("E28").Value = "  +1.70%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  +2.20%  "

# Row 30 - Hedera
Set-TextValue $ws.Range("D30") "0.0514"
$ws.Range("E30").Value = "  -1.27%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.19%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +4.16%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D33") "2.98"
$ws.Range("E33").Value = "  +1.22%  "

# Row 34 - Maker
Set-TextValue $ws.Range("D34") "1.264.59"
$ws.Range("E34").Value = "  +8.26%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +0.86%  "

# Row 36 - HuobiToken
Set-TextValue $ws.Range("D36") "2.38"
$ws.Range("E36").Value = "  +0.16%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  +4.58%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  +1.56%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  -1.01%  "

# Row 41 - MXToken
$ws.Range("E41").Value = "  -1.68%  "

# Row 42 - TrustWalletToken
Set-TextValue $ws.Range("D42") "0.798"
$ws.Range("E42").Value = "  +0.82%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  -0.38%  "

# Row 44 - RocketPoolETH
Set-TextValue $ws.Range("D44") "1.770.70"
$ws.Range("E44").Value = "  +0.15%  "

# Row 45 - Quant
$ws.Range("E45").Value = "  +1.58%  "

# Row 46 - RenderToken
$ws.Range("E46").Value = "  +3.14%  "

# Row 47 - Aave
Set-TextValue $ws.Range("D47") "55.17"
$ws.Range("E47").Value = "  +0.92%  "

# Row 48 - BabyDogeCoin
$ws.Range("E48").Value = "  -2.05%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  +0.22%  "

# Row 50 - EnergySwap
Set-TextValue $ws.Range("D50") "7.56"
$ws.Range("E50").Value = "  +0.16%  "

# Row 51 - now Algorand (was Mantle)
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.0959"
$ws.Range("E51").Value = "  +2.34%  "
